# Update the data table (rows 16-27) on Hoja1 so it is sorted by Periodo Mora
# (ascending, 1901..1906) with each worker's row for that period listed
# together (JORGE LUIS MIRANDA RAMOS then LINA MARIA MAZO MONSALVE), matching
# the refreshed "Estado de Cuenta" database export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$jorge = "73270932"
$jorgeName = "JORGE LUIS MIRANDA RAMOS"
$lina = "43600251"
$linaName = "LINA MARIA MAZO MONSALVE"

$rows = @(
    @{ Row = 16; Doc = $jorge; Name = $jorgeName; Periodo = "1901"; Mora = 40000;  Salario = 1000000 },
    @{ Row = 17; Doc = $lina;  Name = $linaName;  Periodo = "1901"; Mora = 100000; Salario = 2500000 },
    @{ Row = 18; Doc = $jorge; Name = $jorgeName; Periodo = "1902"; Mora = 40000;  Salario = 1000000 },
    @{ Row = 19; Doc = $lina;  Name = $linaName;  Periodo = "1902"; Mora = 100000; Salario = 2500000 },
    @{ Row = 20; Doc = $jorge; Name = $jorgeName; Periodo = "1903"; Mora = 40000;  Salario = 1000000 },
    @{ Row = 21; Doc = $lina;  Name = $linaName;  Periodo = "1903"; Mora = 100000; Salario = 2500000 },
    @{ Row = 22; Doc = $jorge; Name = $jorgeName; Periodo = "1904"; Mora = 40000;  Salario = 1000000 },
    @{ Row = 23; Doc = $lina;  Name = $linaName;  Periodo = "1904"; Mora = 100000; Salario = 2500000 },
    @{ Row = 24; Doc = $jorge; Name = $jorgeName; Periodo = "1905"; Mora = 40000;  Salario = 1000000 },
    @{ Row = 25; Doc = $lina;  Name = $linaName;  Periodo = "1905"; Mora = 100000; Salario = 2500000 },
    @{ Row = 26; Doc = $jorge; Name = $jorgeName; Periodo = "1906"; Mora = 22666;  Salario = 1000000 },
    @{ Row = 27; Doc = $lina;  Name = $linaName;  Periodo = "1906"; Mora = 56666;  Salario = 2500000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 2).Value = "CC"
    $ws.Cells.Item($n, 3).Value = $r.Doc
    $ws.Cells.Item($n, 4).Value = $r.Name
    $ws.Cells.Item($n, 5).Value = $r.Periodo
    $ws.Cells.Item($n, 6).Value = $r.Mora
    $ws.Cells.Item($n, 7).Value = $r.Salario
}
